# Update automatico via Actualizar 06-14-2020 00-51-48
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trabajo")

$newRow = 22

$ws.Range("E$newRow").Value = "https://www.mitradel.gob.pa/verifican-reactivacion-de-contratos-laborales-en-comercios-de-veraguas/"
$ws.Hyperlinks.Add($ws.Range("E$newRow"), "https://www.mitradel.gob.pa/verifican-reactivacion-de-contratos-laborales-en-comercios-de-veraguas/") | Out-Null

$ws.Range("F$newRow").Value = "Inspectores de la Dirección Regional del Ministerio de Trabajo y Desarrollo Laboral (Mitradel) en la provincia de Veraguas realizaron una serie de operativos para verificar la reactivación de los contratos de trabajo en comercios de los bloques 1 y 2. Para poder hacer efectivo este proceso las empresas deben completar el formulario digital de “Reactivación de Contratos”, disponible en la página web www.mitradel.gob.pa."

$ws.Range("H$newRow").Value = "13/08"

$ws.Range("A$newRow").Style = $ws.Range("A21").Style
$ws.Range("B$newRow").Style = $ws.Range("B21").Style
$ws.Range("C$newRow").Style = $ws.Range("C21").Style
$ws.Range("D$newRow").Style = $ws.Range("D21").Style
$ws.Range("E$newRow").Style = $ws.Range("E21").Style
$ws.Range("F$newRow").Style = $ws.Range("F21").Style
$ws.Range("I$newRow").Style = $ws.Range("I21").Style
$ws.Range("J$newRow").Style = $ws.Range("J21").Style
$ws.Range("K$newRow").Style = $ws.Range("K21").Style

$ws.Rows($newRow).RowHeight = 89.25
